# EQR_cashflow.xlsx update
# - Fills in previously-blank "B" (trailing-twelve-month / most-recent-period)
#   column figures for most line items on the cash-flow statement.
# - Corrects a couple of previously-entered figures (B24, F11, F26).
# - Widens column B so it matches the width used by the rest of the data
#   columns (C:AO).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Newly populated column B values (were blank inline strings before).
$ws.Range("B3").Value  = 866033000
$ws.Range("B4").Value  = -297532000
$ws.Range("B5").Value  = -13394000
$ws.Range("B6").Value  = -16267000
$ws.Range("B7").Value  = 133000
$ws.Range("B8").Value  = 1222443000
$ws.Range("B9").Value  = 311494000
$ws.Range("B11").Value = -3468000
$ws.Range("B12").Value = -83851000
$ws.Range("B13").Value = 278004000
$ws.Range("B14").Value = -607112900
$ws.Range("B15").Value = 29376000
$ws.Range("B17").Value = -900097000
$ws.Range("B18").Value = -65547000
$ws.Range("B19").Value = -1543381000
$ws.Range("B20").Value = -42934000
$ws.Range("B21").Value = 99728000
$ws.Range("B22").Value = 56794000
$ws.Range("B23").Value = 24008000
$ws.Range("B25").Value = -29690000
$ws.Range("B26").Value = -3468000
$ws.Range("B27").Value = -3090000
$ws.Range("B28").Value = 29376000

# Corrected figures.
$ws.Range("B24").Value = -897007000
$ws.Range("F11").Value = -8799000
$ws.Range("F26").Value = -8799000

# Match column B's width to the rest of the data columns.
$ws.Columns.Item(2).ColumnWidth = $ws.Columns.Item(3).ColumnWidth
